# Weekly crime-stats refresh: new week's numbers layered onto the
# cs-en-us-018pct workbook (new crime data collected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: bump the volume/issue number and roll the reporting week
# forward by one week. Edit the text in place (via Characters) instead of
# clobbering the whole cell, so only the trailing token changes.
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "35"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 9).Text = "8/28/2023"
$weekCell.Characters(47, 9).Text = "9/3/2023"

# ---------------------------------------------------------------------
# Helper: turn a numeric cell into the workbook's "no data" marker while
# keeping the same look-and-feel (font/alignment/number format) as the
# rest of its column. Count-style columns show a literal "0", percent
# -change columns show "***.*".
# ---------------------------------------------------------------------
function Set-Blank($cell, $refCell, $text) {
    $ws.Range($refCell).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $text
    $ws.Range($refCell).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
}

# Helper: turn a "no data" marker cell back into a real number, picking
# up the normal numeric styling from a same-column reference cell first.
function Set-FromBlank($cell, $refCell, $value) {
    $ws.Range($refCell).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).Value = $value
}

# ---------------------------------------------------------------------
# Row 15 - Precinct 22
# ---------------------------------------------------------------------
Set-Blank "F15" "F14" "0"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
$ws.Range("M15").Value = -41.666666666666

# ---------------------------------------------------------------------
# Row 16 - Precinct 23
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -15.384615384615
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = -17
$ws.Range("L16").Value = 1.219512195121
$ws.Range("M16").Value = -6.741573033707
$ws.Range("N16").Value = -91.065662002152

# ---------------------------------------------------------------------
# Row 17 - Precinct 24
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -85.714285714285
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -54.545454545454
$ws.Range("I17").Value = 116
$ws.Range("J17").Value = 124
$ws.Range("K17").Value = -6.451612903225
$ws.Range("L17").Value = 6.422018348623
$ws.Range("M17").Value = 19.587628865979
$ws.Range("N17").Value = -65.882352941176

# ---------------------------------------------------------------------
# Row 18 - Precinct 25
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -31.25
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 153
$ws.Range("K18").Value = -42.483660130719
$ws.Range("L18").Value = -37.142857142857
$ws.Range("M18").Value = -31.782945736434
$ws.Range("N18").Value = -93.901593901593

# ---------------------------------------------------------------------
# Row 19 - Precinct 26
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 45
$ws.Range("E19").Value = -35.555555555555
$ws.Range("G19").Value = 155
$ws.Range("H19").Value = -16.129032258064
$ws.Range("I19").Value = 1206
$ws.Range("J19").Value = 1268
$ws.Range("K19").Value = -4.889589905362
$ws.Range("L19").Value = 69.620253164557
$ws.Range("M19").Value = 11.977715877437
$ws.Range("N19").Value = -76.13299030279

# ---------------------------------------------------------------------
# Row 20 - Precinct 27 (D/E flip from "no data" to real numbers)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 4
Set-FromBlank "D20" "D19" 1
Set-FromBlank "E20" "E19" 300
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = -30.864197530864
$ws.Range("L20").Value = 55.555555555555
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = -83.333333333333

# ---------------------------------------------------------------------
# Row 21 - Precinct 28 (bold totals row)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -36.206896551724
$ws.Range("F21").Value = 171
$ws.Range("G21").Value = 216
$ws.Range("H21").Value = -20.833333333333
$ws.Range("I21").Value = 1558
$ws.Range("J21").Value = 1741
$ws.Range("K21").Value = -10.511200459506
$ws.Range("L21").Value = 42.413162705667
$ws.Range("M21").Value = 8.044382801664
$ws.Range("N21").Value = -80.836408364083

# ---------------------------------------------------------------------
# Row 22 - Precinct 29 (C/D/E flip from real numbers to "no data")
# ---------------------------------------------------------------------
Set-Blank "C22" "C14" "0"
Set-Blank "D22" "D14" "0"
Set-Blank "E22" "E14" "***.*"
$ws.Range("I22").Value = 47
$ws.Range("K22").Value = 6.818181818181
$ws.Range("L22").Value = 30.555555555555
$ws.Range("M22").Value = 20.512820512820

# ---------------------------------------------------------------------
# Row 23 - Precinct 30 (L flips from "no data" to a real number)
# ---------------------------------------------------------------------
Set-FromBlank "L23" "L19" 200

# ---------------------------------------------------------------------
# Row 24 - Precinct 31
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 54
$ws.Range("D24").Value = 76
$ws.Range("E24").Value = -28.947368421052
$ws.Range("F24").Value = 222
$ws.Range("G24").Value = 240
$ws.Range("H24").Value = -7.5
$ws.Range("I24").Value = 1807
$ws.Range("J24").Value = 1727
$ws.Range("K24").Value = 4.632310364794
$ws.Range("L24").Value = 47.993447993448
$ws.Range("M24").Value = 44.213886671987

# ---------------------------------------------------------------------
# Row 25 - Precinct 32
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 36.842105263157
$ws.Range("I25").Value = 469
$ws.Range("J25").Value = 372
$ws.Range("K25").Value = 26.075268817204
$ws.Range("L25").Value = 66.903914590747
$ws.Range("M25").Value = 53.770491803278

# ---------------------------------------------------------------------
# Row 26 - Precinct 33
# ---------------------------------------------------------------------
Set-Blank "F26" "F14" "0"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -100

# ---------------------------------------------------------------------
# Row 27 - Precinct 34
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 44.444444444444
$ws.Range("I27").Value = 74
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = 13.846153846153
$ws.Range("L27").Value = 25.423728813559

# ---------------------------------------------------------------------
# Row 30 - Precinct 37 (D/E flip from real numbers to "no data")
# ---------------------------------------------------------------------
Set-Blank "D30" "D14" "0"
Set-Blank "E30" "E14" "***.*"
